$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 408 (pushes existing rows 408-524 down to 409-525,
# growing the used range from A1:R524 to A1:R525, matching the dimension change in the diff).
$ws.Rows.Item(408).Insert()

# Populate the newly inserted row 408 with the new data record.
$ws.Range("A408").Value = 10
$ws.Range("B408").Value = 'Vega Modelo de Temuco'
$ws.Range("C408").Value = 'La Araucanía'
$ws.Range("D408").Value = 44841
$ws.Range("E408").Value = 9
$ws.Range("F408").Value = 100112023
$ws.Range("G408").Value = 'Brócoli'
$ws.Range("H408").Value = 'Sin especificar'
$ws.Range("I408").Value = 'Primera'
$ws.Range("J408").Value = 1400
$ws.Range("K408").Value = 1200
$ws.Range("L408").Value = 1200
$ws.Range("M408").Value = 1200
$ws.Range("N408").Value = '$/unidad'
$ws.Range("O408").Value = 'Región Metropolitana'
$ws.Range("P408").Value = 1200
$ws.Range("Q408").Value = 1
$ws.Range("R408").Value = 'Hortaliza'
